# FICHAMENTO.xlsx - "Criacao de topicos contendo artigos relevantes para
# cada topico - Continuidade de leitura"
#
# Adds the "Tema" (column G) classification "Wild Cards" to the rows that
# were still missing it (rows 13-33), creates a new "Weak Signals
# Management" topic for the existing Muniz/Blanck (2014) reference (row 34,
# column G), refreshes that reference's bibkey, and appends a brand new
# row (35) with another fichamento entry for the same "Weak Signals
# Management" topic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 13-33: tag every remaining reading with the "Wild Cards" topic ---
# Rows 13/14 already own an (empty) G cell with the "text" number format;
# rows 15-33 have none at all. Copying the format from G12 (which already
# carries the correct "Wild Cards" tag) keeps every cell on the existing
# style instead of minting new ones, then the values are written in bulk.
$ws.Range("G12").Copy() | Out-Null
$ws.Range("G13:G33").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = 0
$ws.Range("G13:G33").Value = "Wild Cards"

# --- Row 34: existing Muniz/Blanck (2014) reference gets its own topic ---
# Bibkey correction first ("inproceedings" -> "2014MunizMery"), then the new
# topic tag, so the shared-string table grows in the same order as the
# original author's edit.
$bibkey = @"
@inproceedings{2014MunizMery,
author = {Muniz, Raquel and Blanck, Mery},
year = {2014},
month = {05},
pages = {},
title = {WEAK SIGNALS MANAGEMENT, ENTREPRENEURSHIP AND UNCERTAINTY: A RELATIONAL THEORETICAL ESSAY UNDER THE PERSPECTIVE OF INTELLIGENCE},
isbn = {978-859969310-0},
doi = {10.5748/9788599693100-11CONTECSI/COMM-625}
}
"@
$ws.Range("I34").Value = $bibkey

$ws.Range("G34").Value = "Weak Signals Management"

# --- Row 35: new fichamento entry, same "Weak Signals Management" topic ---
$ws.Range("C35").Value = "8"
$ws.Range("D35").Value = "Sobretudo, os aspectos relativos à informação são importantes na medida em que se constata que é exatamente a partir dos dados, do significado da informação que resulta destes e do conhecimento advindo a partir da compreensão, entendimento e aprendizado proporcionados pelo conhecimento que a Inteligência pode se estabelecer."
$ws.Range("G35").Value = "Weak Signals Management"
$ws.Range("H35").Value = "Sobre informações e conhecimento"

# --- View bookkeeping to mirror the saved workbook state ---
$ws.Application.Windows.Item(1).Zoom = 107
$ws.Range("D36").Select() | Out-Null
